# StructureDefinition-employee-shift.xlsx update:
# - bump URL/version/date/publisher metadata on the "Metadata" sheet
# - keep the "Fixed Value" for Extension.url on the "Elements" sheet in sync
#   with the new canonical URL
# - the top-level "Extension" row's rolled-up Constraint(s) text is cleared
#   (it now only lives on the Extension.extension row)

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-shift"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-shift"
$elements.Range("AI2").Value = ""
